$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-13 are updated from 45184 to 45185
# (the date value is incremented by one day, e.g. 2023-09-15 -> 2023-09-16)
for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45185
    }
}
